# "1. added msg send panel"
# - Adds a new "is_running_in_bat" field to bms_cfg (msg: MSG_BMS_GET_RUNNING_IN_BAT)
#   and re-sorts that sheet's fields by c_datatype size (uint8_t, int16_t, uint16_t, float).
# - Adds a new "strokes_count" field to mainpump_cfg (msg: MSG_MAINPUMP_GET_STROKES_COUNT)
#   and re-sorts that sheet's fields (bool first, then uint8_t, then the new float).
# - Fixes a c_datatype mismatch on ccpd_post_st (mainpump_st / end_time_ms swap
#   uint32_t <-> uint8_t).
# - Leaves the active sheet / selection on mainpump_cfg (the newly edited panel).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. bms_cfg: insert "is_running_in_bat" after "is_charging" and re-sort the
#    field table by c_datatype (uint8_t -> int16_t -> uint16_t -> float).
# ---------------------------------------------------------------------------
$wsBms = $wb.Worksheets.Item("bms_cfg")

$bmsRows = @(
    @{D="is_charging";        E="uint8_t";  F=1; G="MSG_BMS_GET_CHARGING_STATUS"},
    @{D="is_running_in_bat";  E="uint8_t";  F=1; G="MSG_BMS_GET_RUNNING_IN_BAT"},
    @{D="rel_soc";            E="uint8_t";  F=1; G="MSG_BMS_GET_RELATIVE_SOC"},
    @{D="abs_soc";            E="uint8_t";  F=1; G="MSG_BMS_GET_ABSOLUTE_SOC"},
    @{D="health_state";       E="uint8_t";  F=1; G="MSG_BMS_GET_BAT_HEALTH"},
    @{D="bat_current";        E="int16_t";  F=2; G="MSG_BMS_GET_BAT_CURRENT"},
    @{D="bat_avg_current";    E="int16_t";  F=2; G="MSG_BMS_GET_BAT_AVG_CRNT"},
    @{D="rem_capacity";       E="uint16_t"; F=2; G="MSG_BMS_GET_REM_CAPACITY"},
    @{D="run_time_to_empty";  E="uint16_t"; F=2; G="MSG_BMS_GET_RUN_TIME_TO_EMPTY"},
    @{D="avg_time_to_empty";  E="uint16_t"; F=2; G="MSG_BMS_GET_AVG_TIME_TO_EMPTY"},
    @{D="avg_time_to_full";   E="uint16_t"; F=2; G="MSG_BMS_GET_AVG_TIME_TO_FULL"},
    @{D="charging_current";   E="uint16_t"; F=2; G="MSG_BMS_GET_BAT_CHG_CURRENT"},
    @{D="num_dschg_cycles";   E="uint16_t"; F=2; G="MSG_BMS_GET_BAT_DSCHG_COUNT"},
    @{D="run_interval_ms";    E="uint16_t"; F=2; G=""},
    @{D="bat_temp";           E="float";    F=4; G="MSG_BMS_GET_BAT_TEMP"},
    @{D="bat_volt";           E="float";    F=4; G="MSG_BMS_GET_BAT_VOLT"},
    @{D="charging_voltage";   E="float";    F=4; G="MSG_BMS_GET_BAT_CHG_VOLT"},
    @{D="cell1_voltage";      E="float";    F=4; G="MSG_BMS_GET_BAT_CELL1_VOLT"},
    @{D="cell2_voltage";      E="float";    F=4; G="MSG_BMS_GET_BAT_CELL2_VOLT"},
    @{D="cell3_voltage";      E="float";    F=4; G="MSG_BMS_GET_BAT_CELL3_VOLT"}
)

$r = 2
foreach ($row in $bmsRows) {
    $wsBms.Cells.Item($r, 4).Value = $row.D
    $wsBms.Cells.Item($r, 5).Value = $row.E
    $wsBms.Cells.Item($r, 6).Value = $row.F
    $wsBms.Cells.Item($r, 7).Value = $row.G
    $r++
}

# ---------------------------------------------------------------------------
# 2. mainpump_cfg: re-sort existing fields and append "strokes_count".
# ---------------------------------------------------------------------------
$wsPump = $wb.Worksheets.Item("mainpump_cfg")

$pumpRows = @(
    @{D="is_pump_started";          E="bool";    F=1; G="MSG_MAINPUMP_GET_DRIVE_STAGE"},
    @{D="cur_pwm_dc";                E="uint8_t"; F=1; G="MSG_MAINPUMP_GET_CUR_PWM_DC"},
    @{D="pwm_dc_incdec_step_size";   E="uint8_t"; F=1; G="MSG_MAINPUMP_GET_PWM_INC_DEC_STEP"},
    @{D="strokes_count";             E="float";   F=4; G="MSG_MAINPUMP_GET_STROKES_COUNT"}
)

$r = 2
foreach ($row in $pumpRows) {
    $wsPump.Cells.Item($r, 4).Value = $row.D
    $wsPump.Cells.Item($r, 5).Value = $row.E
    $wsPump.Cells.Item($r, 6).Value = $row.F
    $wsPump.Cells.Item($r, 7).Value = $row.G
    $r++
}

# ---------------------------------------------------------------------------
# 3. ccpd_post_st: fix mainpump_st / end_time_ms c_datatype swap.
# ---------------------------------------------------------------------------
$wsCcpd = $wb.Worksheets.Item("ccpd_post_st")
$wsCcpd.Range("E3").Value = "uint8_t"
$wsCcpd.Range("E17").Value = "uint32_t"

# ---------------------------------------------------------------------------
# 4. Update selections on the touched sheets, then finish with mainpump_cfg
#    active/selected (matches the saved workbook view).
# ---------------------------------------------------------------------------
$wsCcpd.Range("E17").Select() | Out-Null

$wsUc = $wb.Worksheets.Item("uc_running_cfg")
$wsUc.Range("D3:F3").Select() | Out-Null

$wsBms.Range("C3").Select() | Out-Null

$wsPump.Range("D5").Select() | Out-Null
$wsPump.Activate() | Out-Null
